$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "2. Data reporter" section (B6:B10) with the new contact
# information for the National Statistical Committee of the Kyrgyz Republic.
$ws.Range("B6").Value = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value = "Kalymbetova Yryskan"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Move the active selection from B2 to B6, matching the saved view state.
$ws.Range("B6").Select()

# Update the workbook window position/size recorded in the saved view
# (maximized-like full-screen window starting at the top-left corner).
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 11835
